# Append the 27 March 2020 top-level (total cases only) DHB case numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = 43917

$data = @(
    @("Auckland", 58),
    @("Bay of Plenty", 6),
    @("Canterbury", 30),
    @("Capital and Coast", 42),
    @("Counties Manukau", 28),
    @("Hawke's Bay", 10),
    @("Hutt Valley", 11),
    @("Lakes", 8),
    @("MidCentral", 7),
    @("Nelson Marlborough", 18),
    @("Northland", 4),
    @("South Canterbury", 2),
    @("Southern", 39),
    @("Tairāwhiti", 1),
    @("Taranaki", 6),
    @("Waikato", 45),
    @("Wairarapa", 5),
    @("Waitemata", 47),
    @("West Coast", 1),
    @("Whanganui", 0),
    @("Total", 368)
)

$startRow = 65
$row = $startRow
foreach ($entry in $data) {
    $dhb = $entry[0]
    $count = $entry[1]

    $ws.Cells.Item($row, 1).Value = $dhb
    $ws.Cells.Item($row, 2).Value = "Total cases"
    $ws.Cells.Item($row, 3).Value = $count
    $ws.Cells.Item($row, 4).Value = $date
    $ws.Cells.Item($row, 4).NumberFormat = "yyyy-mm-dd"

    $row = $row + 1
}
